$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.470.34"
$ws.Range("D3").Value = "2.643.29"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.06"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.19"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").Value = "2.672.72"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.44"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D14").Value = "3.102.43"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "59.415.33"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.30"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "2.663.07"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.62"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.15"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.53"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.56"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "2.760.35"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").Value = "0.0₃0820"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.51"
$ws.Range("E32").Value = "  +9.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.13"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.86"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +13.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.08"
$ws.Range("E37").Value = "  +3.39%  "
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.870"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.76"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "287.94"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.619"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0990"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.76"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0545"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.76"
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("E51").Value = "  -1.30%  "
